$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Date format update (numFmt YYYY-MM-DD -> YYYY-MM-DD HH:MM:SS) ---
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Row 14: new total row (label added first so it lands in the shared-string table
#     right after the existing labels, matching the order the strings were authored) ---
$ws.Range("D14").Value = "total ="

# --- Row 2: Pfizer ---
$ws.Range("A2").Value = 45656
$ws.Range("D2").Value = 5340
$ws.Range("E2").Value = 267000
$ws.Range("F2").Value = "Dec 30, 3:57:35 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 3: Ultratech Cement ---
$ws.Range("D3").Value = 11350
$ws.Range("E3").Value = 170250
$ws.Range("F3").Value = "Dec 30, 3:58:44 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 4: Reliance Industries ---
$ws.Range("D4").Value = 1211.5
$ws.Range("E4").Value = 174456
$ws.Range("F4").Value = "Dec 30, 3:59:56 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 5: Aditya Birla Capital ---
$ws.Range("D5").Value = 180
$ws.Range("E5").Value = 23940
$ws.Range("F5").Value = "Dec 30, 3:59:25 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 6: Grasim Industries ---
$ws.Range("D6").Value = 2445
$ws.Range("E6").Value = 232275
$ws.Range("F6").Value = "Dec 30, 3:59:12 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 7: Bombay Dyeing ---
$ws.Range("D7").Value = 185.3
$ws.Range("E7").Value = 11118
$ws.Range("F7").Value = "Dec 30, 3:59:35 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 8: Raymond ---
$ws.Range("D8").Value = 1689.6
$ws.Range("E8").Value = 47308.8
$ws.Range("F8").Value = "Dec 30, 3:58:27 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 9: LICI ---
$ws.Range("D9").Value = 912.5
$ws.Range("E9").Value = 13687.5
$ws.Range("F9").Value = "Dec 30, 3:59:58 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 10: Tata Steel ---
$ws.Range("D10").Value = 137.5
$ws.Range("E10").Value = 46062.5
$ws.Range("F10").Value = "Dec 30, 3:59:57 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 11: Mahindra Mahindra ---
$ws.Range("D11").Value = 3015
$ws.Range("E11").Value = 144720
$ws.Range("F11").Value = "Dec 30, 3:59:29 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 12: Jio Financial ---
$ws.Range("D12").Value = 304.8
$ws.Range("E12").Value = 21945.6
$ws.Range("F12").Value = "Dec 30, 3:59:57 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 13: Raymond Lifestyle ---
$ws.Range("D13").Value = 2115.5
$ws.Range("E13").Value = 46541
$ws.Range("F13").Value = "Dec 30, 3:56:51 PM GMT+5:30 · INR · NSE · Disclaimer"

# --- Row 14: new total row (value) ---
$ws.Range("E14").Value = 1199304.4
